$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Usuarios")
$ws2 = $wb.Worksheets.Item("Recursos")

# ---------------------------------------------------------------------------
# Sheet "Usuarios": reorder the last four names and move the selection
# ---------------------------------------------------------------------------
$ws1.Range("A15").Value = "Sanchez Barreiro, Pablo"
$ws1.Range("A16").Value = "Pedro"
$ws1.Range("A17").Value = "JAVI"
$ws1.Range("A18").Value = "RODRIGUEZ PÉREZ, DANIEL"

$ws1.Range("A16").Select()

$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Sheet "Recursos": add an "IDRecurso" column, an "Excluido" column, and a
# new row, then fill in the resource list (identified by ID now).
# ---------------------------------------------------------------------------

# Make room for the new "IDRecurso" column between "Contexto del evento"
# and "Alias" (old column B becomes column C).
$ws2.Columns.Item(2).Insert()

$ws2.Range("B1").Value = "IDRecurso"

# Add the trailing "Excluido" header, copying the header style from C1 so
# it keeps the same bold/border/alignment formatting.
$ws2.Range("C1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws2.Range("D1").Value = "Excluido"

# Resource rows: context label, resource id, alias (defaults to the label).
$data = @(
    @("Foro: Noticias de clase", 5000),
    @("Carpeta: Exámenes", 5002),
    @("Carpeta: Recursos del Alumnado", 5011),
    @("Carpeta: Recursos del Alumnado", 5012),
    @("Carpeta: Papeleo", 5013),
    @("Carpeta: Recursos del Alumnado", 5014),
    @("Tarea: Entrega inicial", 5015),
    @("Carpeta: Entrega inicial", 5016)
)

$row = 2
foreach ($item in $data) {
    $ws2.Range("A$row").Value = $item[0]
    $ws2.Range("B$row").Value = $item[1]
    $ws2.Range("C$row").Value = $item[0]
    $row = $row + 1
}

# Column widths: narrow hidden ID column, equal-width label/alias columns.
$ws2.Columns.Item(1).ColumnWidth = 29.8
$ws2.Columns.Item(2).ColumnWidth = -0.9
$ws2.Columns.Item(2).Hidden = $true
$ws2.Columns.Item(3).ColumnWidth = 29.8
